$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format on cells whose new numeric-looking value
# would otherwise be auto-converted to a genuine Number by Excel (losing the
# original text formatting, e.g. trailing zeros or plain decimal notation).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D10",
    "D11",
    "D14",
    "D20",
    "D24",
    "D25",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D39",
    "D42",
    "D44",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values exactly as captured from the source refresh.
$ws.Range("D2").Value = "63.913.77"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.080.30"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "537.07"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "135.91"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.076.08"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  +5.17%  "
$ws.Range("D14").Value = "34.36"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "3.573.59"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "63.902.07"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "3.080.84"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "483.28"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "79.71"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").Value = "12.12"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "26.20"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("D32").Value = "1.88"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "57.09"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").Value = "2.34"
$ws.Range("E34").Value = "  -6.35%  "
$ws.Range("D35").Value = "499.77"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").Value = "5.36"
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "3.234.19"
$ws.Range("E38").Value = "  +6.41%  "
$ws.Range("D39").Value = "0.0396"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "8.09"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "0.254"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "2.03"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "121.74"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "0.0₃0527"
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "24.31"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "2.40"
$ws.Range("E51").Value = "  +3.56%  "
